# Generate Report for Handoff
#
# Status moves from "In Translation" to "Ready for handoff", and the
# "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" timestamps
# advance a few seconds to record the handoff moment. The Status columns
# were also widened (to fit the new, longer status text) on every sheet
# that shows it.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-04 16:43:52"
# Widen the zh-cn / de-de status columns (E:F) to fit the longer text.
# ColumnWidth is quantized to whole pixels by this host, so feed the
# chars value whose nearest on-grid pixel width best matches the target.
$wsOverview.Range("E2:F2").ColumnWidth = 16.333333333333332

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-04 16:43:48"
$wsZhCn.Range("C2").ColumnWidth = 16.333333333333332

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-04 16:43:52"
$wsDeDe.Range("C2").ColumnWidth = 16.333333333333332
